$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Row 5: swap the 2-pin terminal block supplier from TE/Digi-key to
# DBParts/Amazon. Name (A5) and Description (B5) both become the new
# product's name; Vendor (F5) becomes "Amazon"; the link (G5) now
# resolves to the Amazon URL and is styled like the other link cells.
# -----------------------------------------------------------------
$ws.Range("A5").Value = "DBParts 20 pc 2-pin 0.1” pitch Terminal Blocks"
$ws.Range("B5").Value = "DBParts 20 pc 2-pin 0.1” pitch Terminal Blocks"
$ws.Range("C5").Value = 6.99
$ws.Range("D5").Formula = "=1/20"
$ws.Range("F5").Value = "Amazon"
$ws.Range("G5").Value = "https://www.amazon.com/DBParts-20pcs-Terminal-Connector-2-54mm/dp/B07NSJV6NW/ref=sxbs_sxwds-stvp?cv_ct_cx=terminal+block+assortment&keywords=terminal+block+assortment&pd_rd_i=B07NSJV6NW&pd_rd_r=f9187fb7-4ab9-45b4-a445-3bf2b68a1d13&pd_rd_w=WAjjn&pd_rd_wg=olYB9&pf_rd_p=a6d018ad-f20b-46c9-8920-433972c7d9b7&pf_rd_r=PMA44C3EHR468DBRJFP7&psc=1&qid=1581446551&sr=1-3-dd5817a1-1ba7-46c2-8996-f96e7b0f409c"
$ws.Range("G5").Font.Color = 16711680

# -----------------------------------------------------------------
# Row 11: swap the 7-segment display supplier from Lite-On/Digi-key
# to Uxcell/Amazon. Vendor (F11) becomes "Amazon"; quantity becomes a
# fraction (9 boards per pack of 10); price and link update.
# -----------------------------------------------------------------
$ws.Range("A11").Value = "Uxcell 10 pc 7-Segment 10DIN LED Display"
$ws.Range("B11").Value = "7-Segment 10-DIN 10 pc"
$ws.Range("C11").Value = 5.99
$ws.Range("D11").Formula = "=9/10"
$ws.Range("F11").Value = "Amazon"
$ws.Range("G11").Value = "https://www.amazon.com/a13071500ux0900-Cathode-Segment-Display-Digital/dp/B00EZBGUMC/ref=sr_1_fkmr0_1?keywords=7-Segment+10+DIN+LED&qid=1581447942&sr=8-1-fkmr0"

# -----------------------------------------------------------------
# Cosmetic: widen the Name/Description columns and move the
# selection cursor to the new C14:C15 block.
# (ColumnWidth is quantized to pixel boundaries by the host the same
# way real Excel quantizes character widths, so the inputs below are
# chosen to land on the nearest achievable width to the 40.01 / 39.32
# targets after that rounding.)
# -----------------------------------------------------------------
$ws.Range("A1").EntireColumn.ColumnWidth = 39.15
$ws.Range("B1").EntireColumn.ColumnWidth = 38.5

$ws.Range("C14:C15").Select()
